$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These D-column cells get new values that Excel's automatic type
# detection would parse as numbers (e.g. "234.32"). The source workbook
# stores every Price/Volume cell as text, so force a Text number format
# before writing them, keeping the stored type a string like the original.
$numericLookingRange = $ws.Range("D5,D7,D9,D10,D13,D14,D15,D16,D19,D20,D22,D25,D26,D28,D29,D33,D34,D35,D36,D39,D40,D41,D43,D46,D49")
$numericLookingRange.NumberFormat = "@"

# Write the refreshed price / volume(1h) values scraped by this run.
$ws.Range('D5').Value = '234.32'
$ws.Range('D7').Value = '58.75'
$ws.Range('D9').Value = '0.391'
$ws.Range('D10').Value = '0.0787'
$ws.Range('D13').Value = '14.84'
$ws.Range('D14').Value = '21.08'
$ws.Range('D15').Value = '0.773'
$ws.Range('D16').Value = '5.31'
$ws.Range('D19').Value = '6.15'
$ws.Range('D20').Value = '71.28'
$ws.Range('D22').Value = '228.73'
$ws.Range('D25').Value = '2.41'
$ws.Range('D26').Value = '169.18'
$ws.Range('D28').Value = '9.00'
$ws.Range('D29').Value = '19.53'
$ws.Range('D33').Value = '0.0632'
$ws.Range('D34').Value = '4.66'
$ws.Range('D35').Value = '2.48'
$ws.Range('D36').Value = '1.83'
$ws.Range('D39').Value = '5.40'
$ws.Range('D40').Value = '0.0980'
$ws.Range('D41').Value = '98.45'
$ws.Range('D43').Value = '2.86'
$ws.Range('D46').Value = '1.16'
$ws.Range('D49').Value = '7.46'
$ws.Range('D2').Value = '37.758.04'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '2.081.47'
$ws.Range('E3').Value = '  -0.32%  '
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('E7').Value = '  +1.04%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +0.51%  '
$ws.Range('E10').Value = '  +0.94%  '
$ws.Range('E11').Value = '  +2.92%  '
$ws.Range('D12').Value = '2.389.85'
$ws.Range('E12').Value = '  -0.27%  '
$ws.Range('E13').Value = '  +2.22%  '
$ws.Range('E14').Value = '  -1.08%  '
$ws.Range('E15').Value = '  -1.98%  '
$ws.Range('E16').Value = '  +1.62%  '
$ws.Range('D17').Value = '2.076.64'
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').Value = '37.691.78'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('E19').Value = '  -0.54%  '
$ws.Range('E20').Value = '  +2.19%  '
$ws.Range('D21').Value = '0.0₃0833'
$ws.Range('E21').Value = '  +1.33%  '
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('E24').Value = '  -1.46%  '
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('E27').Value = '  +3.53%  '
$ws.Range('E28').Value = '  +0.90%  '
$ws.Range('E29').Value = '  +1.39%  '
$ws.Range('E30').Value = '  -2.22%  '
$ws.Range('E31').Value = '  +1.97%  '
$ws.Range('E32').Value = '  +0.96%  '
$ws.Range('E33').Value = '  +1.88%  '
$ws.Range('E34').Value = '  +1.59%  '
$ws.Range('E35').Value = '  -1.95%  '
$ws.Range('E36').Value = '  +2.69%  '
$ws.Range('E37').Value = '  -3.88%  '
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('E39').Value = '  -3.51%  '
$ws.Range('E40').Value = '  +2.37%  '
$ws.Range('E41').Value = '  +1.17%  '
$ws.Range('E42').Value = '  +1.23%  '
$ws.Range('E43').Value = '  -2.60%  '
$ws.Range('D44').Value = '1.461.25'
$ws.Range('E44').Value = '  -1.65%  '
$ws.Range('E45').Value = '  +4.96%  '
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('E47').Value = '  +6.25%  '
$ws.Range('E48').Value = '  +1.89%  '
$ws.Range('E49').Value = '  +2.27%  '
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').Value = '2.275.46'
$ws.Range('E51').Value = '  -0.35%  '

# Put those cells back on the default (Normal) style so only their text
# content changed, matching the rest of the untouched-format column.
$numericLookingRange.Style = "Normal"
